$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metric values for rows 2-6 (Logistic_Regression, Random_Forest,
# Kernel_SVM, CatBoost, DNN) across columns B (Train_Accuracy) through
# L (AUC), leaving column A (Model_Name) untouched.

$data = @(
    @(0.87, 0.77, 11, 35, 5, 9, 0.6899999999999999, 0.55, 0.61, 0.88, 0.71),
    @(0.98, 0.77, 10, 36, 4, 10, 0.71, 0.5, 0.59, 0.9, 0.7),
    @(0.88, 0.78, 11, 36, 4, 9, 0.73, 0.55, 0.63, 0.9, 0.72),
    @(0.95, 0.78, 11, 36, 4, 9, 0.73, 0.55, 0.63, 0.9, 0.72),
    @(0.95, 0.82, 20, 6, 34, 0, 0.37, 1, 0.54, 0.15, 0.57)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value2 = $rowValues[$j]
    }
}
